$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.101.58'
$ws.Cells.Item(2, 5).Value = '  -2.25%  '

$ws.Cells.Item(3, 4).Value = '2.496.86'
$ws.Cells.Item(3, 5).Value = '  -4.05%  '

$ws.Cells.Item(4, 5).Value = '  +0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '555.16'
$ws.Cells.Item(5, 5).Value = '  -3.18%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '147.29'
$ws.Cells.Item(6, 5).Value = '  -4.77%  '

$ws.Cells.Item(7, 5).Value = '  +0.17%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.608'
$ws.Cells.Item(8, 5).Value = '  -2.06%  '

$ws.Cells.Item(9, 4).Value = '2.491.82'
$ws.Cells.Item(9, 5).Value = '  -4.20%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.109'
$ws.Cells.Item(10, 5).Value = '  -7.10%  '

$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.155'
$ws.Cells.Item(11, 5).Value = '  -0.86%  '

$ws.Cells.Item(12, 2).Value = 'Toncoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.43'
$ws.Cells.Item(12, 5).Value = '  -6.45%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.362'
$ws.Cells.Item(13, 5).Value = '  -4.93%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '26.38'
$ws.Cells.Item(14, 5).Value = '  -6.44%  '

$ws.Cells.Item(15, 4).Value = '2.957.16'
$ws.Cells.Item(15, 5).Value = '  -3.88%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000167'
$ws.Cells.Item(16, 5).Value = '  -5.84%  '

$ws.Cells.Item(17, 4).Value = '61.981.60'
$ws.Cells.Item(17, 5).Value = '  -2.20%  '

$ws.Cells.Item(18, 4).Value = '2.508.78'
$ws.Cells.Item(18, 5).Value = '  -4.57%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.29'
$ws.Cells.Item(19, 5).Value = '  -5.63%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.05'
$ws.Cells.Item(20, 5).Value = '  -5.45%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.26'
$ws.Cells.Item(21, 5).Value = '  -5.97%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '325.47'
$ws.Cells.Item(22, 5).Value = '  -4.87%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '64.86'
$ws.Cells.Item(24, 5).Value = '  -3.11%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.75'
$ws.Cells.Item(25, 5).Value = '  -0.53%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.0000104'
$ws.Cells.Item(26, 5).Value = '  -2.67%  '

$ws.Cells.Item(27, 4).Value = '2.645.97'
$ws.Cells.Item(27, 5).Value = '  -3.06%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.52'
$ws.Cells.Item(28, 5).Value = '  -2.20%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.57'
$ws.Cells.Item(29, 5).Value = '  -5.88%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '544.76'
$ws.Cells.Item(30, 5).Value = '  -6.94%  '

$ws.Cells.Item(31, 5).Value = '  -0.05%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.74'
$ws.Cells.Item(32, 5).Value = '  -1.59%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.153'
$ws.Cells.Item(33, 5).Value = '  -4.06%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.92'
$ws.Cells.Item(34, 5).Value = '  -6.15%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.60'
$ws.Cells.Item(35, 5).Value = '  -6.73%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.98'
$ws.Cells.Item(36, 5).Value = '  -8.02%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.89'
$ws.Cells.Item(37, 5).Value = '  -8.68%  '

$ws.Cells.Item(38, 5).Value = '  +0.12%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.381'
$ws.Cells.Item(39, 5).Value = '  -5.41%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '18.76'
$ws.Cells.Item(40, 5).Value = '  -4.34%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '147.67'
$ws.Cells.Item(41, 5).Value = '  -4.59%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.72'
$ws.Cells.Item(42, 5).Value = '  -7.26%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  -0.01%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '40.93'
$ws.Cells.Item(44, 5).Value = '  -1.15%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.35'
$ws.Cells.Item(45, 5).Value = '  -3.64%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '149.41'
$ws.Cells.Item(46, 5).Value = '  -3.86%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.66'
$ws.Cells.Item(47, 5).Value = '  -5.45%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '21.71'
$ws.Cells.Item(48, 5).Value = '  -5.94%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0545'
$ws.Cells.Item(49, 5).Value = '  -6.89%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.595'
$ws.Cells.Item(50, 5).Value = '  -5.16%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0959'
$ws.Cells.Item(51, 5).Value = '  -4.11%  '
